$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $v = $cell.Value2()
    if ($v -eq 45243) {
        $cell.Value = 45244
    }
}
